$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $text) {
    $cell = $ws.Range($rangeAddr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "57.394.43"
$ws.Range("E2").Value = "  -0.89%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.089.45"
$ws.Range("E3").Value = "  +0.71%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.05%  "

# Row 5 - BNB
Set-TextValue "D5" "521.35"
$ws.Range("E5").Value = "  +0.74%  "

# Row 6 - Solana
Set-TextValue "D6" "140.67"
$ws.Range("E6").Value = "  -1.78%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.00%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "3.091.53"
$ws.Range("E8").Value = "  +0.88%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  +0.11%  "

# Row 10 - Toncoin
$ws.Range("E10").Value = "  -1.59%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -0.64%  "

# Row 12 - Cardano
Set-TextValue "D12" "0.383"
$ws.Range("E12").Value = "  +1.60%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "3.622.84"
$ws.Range("E13").Value = "  +0.82%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  +0.99%  "

# Row 15 - Avalanche
Set-TextValue "D15" "25.77"
$ws.Range("E15").Value = "  -2.12%  "

# Row 16 - ShibaInu
Set-TextValue "D16" "0.0000162"
$ws.Range("E16").Value = "  -0.88%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "57.509.64"
$ws.Range("E17").Value = "  -0.61%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.092.11"
$ws.Range("E18").Value = "  +0.89%  "

# Row 19 - Polkadot
Set-TextValue "D19" "6.08"
$ws.Range("E19").Value = "  -0.43%  "

# Row 20 - Chainlink
Set-TextValue "D20" "12.73"
$ws.Range("E20").Value = "  -0.77%  "

# Row 21 - Uniswap
Set-TextValue "D21" "8.03"
$ws.Range("E21").Value = "  -0.94%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "339.13"
$ws.Range("E22").Value = "  +2.02%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.18%  "

# Row 24 - Polygon
$ws.Range("E24").Value = "  +1.75%  "

# Row 25 - Litecoin
Set-TextValue "D25" "66.49"
$ws.Range("E25").Value = "  +1.60%  "

# Row 26 - Kaspa
$ws.Range("E26").Value = "  -1.18%  "

# Row 27 - Binance-PegBSC-USD
$ws.Range("E27").Value = "  +0.13%  "

# Row 28 - PEPE
$subscript3 = [string][char]0x2083
$ws.Range("D28").Value = "0.0" + $subscript3 + "0905"
$ws.Range("E28").Value = "  +0.22%  "

# Row 29 - USDe
$ws.Range("E29").Value = "  +0.00%  "

# Row 30 - RenderToken
Set-TextValue "D30" "6.45"
$ws.Range("E30").Value = "  -0.48%  "

# Row 31 - InternetComputer(DFINITY)
Set-TextValue "D31" "7.14"
$ws.Range("E31").Value = "  -1.26%  "

# Row 32 - PancakeSwap
$ws.Range("E32").Value = "  +1.91%  "

# Row 33 - EthereumClassic
Set-TextValue "D33" "20.84"
$ws.Range("E33").Value = "  +0.66%  "

# Row 34 - Fetch.AI
$ws.Range("E34").Value = "  -1.63%  "

# Row 35 - Monero
Set-TextValue "D35" "156.67"
$ws.Range("E35").Value = "  +1.07%  "

# Row 36 - NEARProtocol
Set-TextValue "D36" "4.59"
$ws.Range("E36").Value = "  +1.16%  "

# Row 37 - Aptos
$ws.Range("E37").Value = "  +1.30%  "

# Row 38 - EnergySwap
Set-TextValue "D38" "26.97"
$ws.Range("E38").Value = "  -0.03%  "

# Row 39 - ImmutableX
$ws.Range("E39").Value = "  +0.24%  "

# Row 40 - Hedera
$ws.Range("E40").Value = "  -3.06%  "

# Row 41 - Stacks
Set-TextValue "D41" "1.52"
$ws.Range("E41").Value = "  +10.88%  "

# Row 42 - was Filecoin, now RenzoRestakedETH
$ws.Range("B42").Value = "RenzoRestakedETH"
$ws.Range("C42").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D42").Value = "3.133.23"
$ws.Range("E42").Value = "  +0.80%  "

# Row 43 - was RenzoRestakedETH, now Filecoin
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D43" "3.92"
$ws.Range("E43").Value = "  -0.14%  "

# Row 44 - Mantle
$ws.Range("E44").Value = "  +4.23%  "

# Row 45 - OKB
Set-TextValue "D45" "36.70"
$ws.Range("E45").Value = "  +0.57%  "

# Row 46 - FirstDigitalUSD
Set-TextValue "D46" "1.00"
$ws.Range("E46").Value = "  +0.07%  "

# Row 47 - Maker
$ws.Range("D47").Value = "2.298.22"
$ws.Range("E47").Value = "  +1.58%  "

# Row 49 - ONDO
$ws.Range("E49").Value = "  +3.28%  "

# Row 50 - InjectiveProtocol
Set-TextValue "D50" "20.46"
$ws.Range("E50").Value = "  -1.80%  "

# Row 51 - Cosmos
$ws.Range("E51").Value = "  +1.29%  "
